$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D22").Value = 90.874795959596
$ws.Range("E22").Value = 4.37635199353916
$ws.Range("F22").Value = 90.64365
$ws.Range("G22").Value = 91.0742275
$ws.Range("I22").Value = 78.5208
$ws.Range("J22").Value = 98.6
$ws.Range("L22").Value = -0.365022989356409
$ws.Range("M22").Value = 0.100208204576704
$ws.Range("N22").Value = 0.311013793348951
$ws.Range("O22").Value = 88.225475
$ws.Range("P22").Value = 90.64365
$ws.Range("Q22").Value = 94.20815
$ws.Range("R22").Value = 98.6
$ws.Range("D25").Value = -2.22863131313131
$ws.Range("E25").Value = 1.45828729917437
$ws.Range("F25").Value = -1.9813
$ws.Range("G25").Value = -2.061786875
$ws.Range("H25").Value = 1.16880771
$ws.Range("I25").Value = -7.983
$ws.Range("J25").Value = -0.0032
$ws.Range("K25").Value = 7.9798
$ws.Range("L25").Value = -1.32785639044366
$ws.Range("M25").Value = 2.43244052505035
$ws.Range("N25").Value = 0.103635965612088
$ws.Range("O25").Value = -2.900325
$ws.Range("P25").Value = -1.9813
$ws.Range("Q25").Value = -1.23915
$ws.Range("R25").Value = -0.0032
$ws.Range("D29").Value = 0.544448484848485
$ws.Range("F29").Value = 0.7589
$ws.Range("G29").Value = 0.651705625
$ws.Range("I29").Value = -2.3837
$ws.Range("J29").Value = 1
$ws.Range("O29").Value = 0.4458
$ws.Range("P29").Value = 0.7589
$ws.Range("Q29").Value = 0.930125
$ws.Range("R29").Value = 1
$ws.Range("D30").Value = -25.004696969697
$ws.Range("F30").Value = -24.93
$ws.Range("G30").Value = -24.9524375
$ws.Range("H30").Value = 0.541149
$ws.Range("I30").Value = -27.2
$ws.Range("J30").Value = -23.4
$ws.Range("L30").Value = -0.884112948353181
$ws.Range("M30").Value = 1.61535932391407
$ws.Range("N30").Value = 0.0490118559299037
$ws.Range("O30").Value = -25.2975
$ws.Range("P30").Value = -24.93
$ws.Range("Q30").Value = -24.585
$ws.Range("R30").Value = -23.4
$ws.Range("D37").Value = -58.0878787878788
$ws.Range("F37").Value = -58.3
$ws.Range("G37").Value = -58.15875
$ws.Range("I37").Value = -80.7
$ws.Range("J37").Value = -31.5
$ws.Range("M37").Value = -0.245783325575439
$ws.Range("O37").Value = -65.9
$ws.Range("P37").Value = -58.3
$ws.Range("Q37").Value = -51.3
$ws.Range("R37").Value = -31.5
$ws.Range("D38").Value = -13.0873388888889
$ws.Range("E38").Value = 5.7586626340977
$ws.Range("F38").Value = -12.4029
$ws.Range("G38").Value = -12.80129125
$ws.Range("I38").Value = -30.5596
$ws.Range("J38").Value = -2.1312
$ws.Range("L38").Value = -0.484393589138627
$ws.Range("M38").Value = -0.0773252414139582
$ws.Range("N38").Value = 0.409250332946639
$ws.Range("O38").Value = -16.245825
$ws.Range("P38").Value = -12.4029
$ws.Range("Q38").Value = -8.806575
$ws.Range("R38").Value = -2.1312
$ws.Range("D43").Value = 25.0201141414141
$ws.Range("F43").Value = 24.8354
$ws.Range("G43").Value = 24.91794125
$ws.Range("I43").Value = 10.6684
$ws.Range("J43").Value = 40.4703
$ws.Range("O43").Value = 20.9044
$ws.Range("P43").Value = 24.8354
$ws.Range("Q43").Value = 29.271325
$ws.Range("R43").Value = 40.4703
$ws.Range("D49").Value = -33.8530303030303
$ws.Range("F49").Value = -33.2
$ws.Range("G49").Value = -33.65875
$ws.Range("H49").Value = 5.18910000000001
$ws.Range("I49").Value = -46.9
$ws.Range("J49").Value = -25.3
$ws.Range("L49").Value = -0.428697808412544
$ws.Range("M49").Value = -0.499300335447264
$ws.Range("O49").Value = -37.4
$ws.Range("P49").Value = -33.2
$ws.Range("Q49").Value = -30.3
$ws.Range("R49").Value = -25.3
$ws.Range("D51").Value = -51.0501090909091
$ws.Range("F51").Value = -50
$ws.Range("G51").Value = -50.028105625
$ws.Range("I51").Value = -90.9091
$ws.Range("J51").Value = -25
$ws.Range("M51").Value = 0.267160685303492
$ws.Range("O51").Value = -58.6051
$ws.Range("P51").Value = -50
$ws.Range("Q51").Value = -41.6667
$ws.Range("R51").Value = -25
$ws.Range("C54").Value = 198
$ws.Range("D54").Value = 80.9227272727273
$ws.Range("E54").Value = 19.521817513277
$ws.Range("G54").Value = 80.985625
$ws.Range("H54").Value = 20.53401
$ws.Range("L54").Value = -0.013125088349585
$ws.Range("M54").Value = -0.353177704531993
$ws.Range("N54").Value = 1.38735515946472
$ws.Range("O54").Value = 66.875
$ws.Range("Q54").Value = 94.45
$ws.Range("D57").Value = -14.697551010101
$ws.Range("E57").Value = 9.69740557584923
$ws.Range("F57").Value = -13.6129
$ws.Range("G57").Value = -14.31862125
$ws.Range("H57").Value = 11.90512974
$ws.Range("I57").Value = -42.3358
$ws.Range("J57").Value = -0.1459
$ws.Range("K57").Value = 42.1899
$ws.Range("L57").Value = -0.297055409708941
$ws.Range("M57").Value = -0.790004512241786
$ws.Range("N57").Value = 0.689164605187315
$ws.Range("O57").Value = -21.94075
$ws.Range("P57").Value = -13.6129
$ws.Range("Q57").Value = -7.118975
$ws.Range("R57").Value = -0.1459
$ws.Range("D58").Value = -108.128758080808
$ws.Range("E58").Value = 7.84400781568312
$ws.Range("F58").Value = -108.4105
$ws.Range("G58").Value = -107.9563725
$ws.Range("H58").Value = 8.55430547999999
$ws.Range("I58").Value = -126.5916
$ws.Range("J58").Value = -90.8852
$ws.Range("K58").Value = 35.7064
$ws.Range("L58").Value = -0.139306013846975
$ws.Range("M58").Value = -0.742058031615904
$ws.Range("N58").Value = 0.557449361800882
$ws.Range("O58").Value = -114.0457
$ws.Range("P58").Value = -108.4105
$ws.Range("Q58").Value = -102.396275
$ws.Range("R58").Value = -90.8852
$ws.Range("C59").Value = 197
$ws.Range("D59").Value = 134.152284263959
$ws.Range("E59").Value = 378.908915818964
$ws.Range("G59").Value = 55.4213836477987
$ws.Range("L59").Value = 5.08656419719666
$ws.Range("M59").Value = 25.2728902908518
$ws.Range("N59").Value = 26.9961424662183
$ws.Range("O59").Value = 21
